$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("createUser")
$ws.Range("A2").Value = 160
